# Apply the "testovaci_scenare" worksheet update:
#  - Clear the stray "Přihlaseni " text left over in E13
#  - Extend the test-case table (D:G) with 7 new rows (TS13..TS19) in the
#    same visual style as the existing TS09..TS12 rows
#  - Fix the E12/F12 cell formatting so it matches the bordered/centered
#    look used by the rest of the table (D12/G12 already had it)
#  - Move the active selection to E4
#
# xlPasteFormats
$xlPasteFormats = -4122

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- 1. Drop the leftover text value in E13 (kept blank like F13/G13) ---
$ws.Range("E13").ClearContents()

# --- 2. Build rows 17-23 (TS13..TS19) reusing the formatting already
#        applied to the existing TS09..TS12 block (D13:G16) ---
$ws.Range("D13:G13").Copy()
$ws.Range("D17:G23").PasteSpecial($xlPasteFormats)

$ws.Range("D17").Value = "TS13"
$ws.Range("D18").Value = "TS14"
$ws.Range("D19").Value = "TS15"
$ws.Range("D20").Value = "TS16"
$ws.Range("D21").Value = "TS17"
$ws.Range("D22").Value = "TS18"
$ws.Range("D23").Value = "TS19"

# --- 3. Make E12/F12 use the same bordered/centered format as D12/G12 ---
$ws.Range("D12").Copy()
$ws.Range("E12:F12").PasteSpecial($xlPasteFormats)

# --- 4. Update the saved selection ---
$ws.Range("E4").Select() | Out-Null
